$d = $word.ActiveDocument

# 1) Move the "_GoBack" bookmark from before the "Fred, " paragraph to a
#    point inside the "Par rapport au document ''" paragraph (between
#    "Par ra" and "pport ..."), computed on the *original* (pre-insert)
#    character offsets. Word only ever keeps a single "_GoBack" bookmark,
#    so adding a new one automatically removes the old bookmarkStart/End.
$pRapport = $d.Paragraphs(6)
$target = $pRapport.Range.Start + 6
$bmRange = $d.Range($target, $target)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 2) Insert a new bold paragraph right after "Fred, " and before
#    "Merci de m'avoir ...".
$pFred = $d.Paragraphs(2)
$pFred.Range.InsertParagraphAfter()
$pNew = $d.Paragraphs(3)
$pNew.Range.Text = "JE RAJOUTE CETTE LIGNE POUR VOIR SI LE CHANGEMENT EST REPERE PAR GIT"
$pNew.Range.Bold = 1
